$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (C) for rows 2..12 from "N" to "Y"
$ws.Range("C2:C12").Value = "Y"
